# Generate Report for Handoff
# Swap the 2e4c8da1-... / 8c4d9bc7-... records between row 2 and row 3 on each
# sheet, and mark the 2e4c8da1-... record (now row 3) as "Ready for handoff"
# with fresh timestamps, reflecting a new localization-status report.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("D2").Value = "2016-03-24 10:59:33"

$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-24 11:01:05"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md", [Type]::Missing, [Type]::Missing, "8c4d9bc7-a25e-4855-a451-08d4b76642a9.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/8c4d9bc7-a25e-4855-a451-08d4b76642a9.md", [Type]::Missing, [Type]::Missing, "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md")

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "In Translation"
$zh.Range("E2").Value = "2016-03-24 10:59:29"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"

$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "2016-03-24 11:00:58"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md", [Type]::Missing, [Type]::Missing, "8c4d9bc7-a25e-4855-a451-08d4b76642a9.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c7c147b301660c6e90b0426f84d499fa52f442f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8c4d9bc7-a25e-4855-a451-08d4b76642a9.164af497f17a2b5f94b17116ebbce596d7e2684c.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/8c4d9bc7-a25e-4855-a451-08d4b76642a9.md", [Type]::Missing, [Type]::Missing, "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c7c147b301660c6e90b0426f84d499fa52f442f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8c4d9bc7-a25e-4855-a451-08d4b76642a9.164af497f17a2b5f94b17116ebbce596d7e2684c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.zh-cn.xlf")

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "In Translation"
$de.Range("E2").Value = "2016-03-24 10:59:33"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"

$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "2016-03-24 11:01:05"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md", [Type]::Missing, [Type]::Missing, "8c4d9bc7-a25e-4855-a451-08d4b76642a9.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/585c0fb05df84e4313d95f7c25891608a12d65fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.de-de.xlf", [Type]::Missing, [Type]::Missing, "8c4d9bc7-a25e-4855-a451-08d4b76642a9.164af497f17a2b5f94b17116ebbce596d7e2684c.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a6f43eab59ae847a84cf2d946f8ca65a116416c0/e2e/8c4d9bc7-a25e-4855-a451-08d4b76642a9.md", [Type]::Missing, [Type]::Missing, "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/585c0fb05df84e4313d95f7c25891608a12d65fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8c4d9bc7-a25e-4855-a451-08d4b76642a9.164af497f17a2b5f94b17116ebbce596d7e2684c.de-de.xlf", [Type]::Missing, [Type]::Missing, "2e4c8da1-f4e6-45fb-800a-1e713d47fe8b.973ddf15c8a25d98ff6766fa1398ed5d996b7f50.de-de.xlf")

$wb.Save()
